$d = $word.ActiveDocument

# The edit rewrites the first paragraph of "Ron's part": it moves the
# paragraph-mark run-properties into a proper <w:pPr>, and it inserts a new
# sentence describing how the hashtags were chosen in between "...we can
# label. " and "We used tweepy ...". A couple of proofErr (spell-check)
# markers are also introduced around "kaggle" and "tweepy".

$p = $d.Paragraphs(1)
$rng = $p.Range

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="774C912F" w14:textId="54ADBA02" w:rsidR="00F318C5" w:rsidRPr="00667113" w:rsidRDefault="00667113" w:rsidP="00667113" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r w:rsidRPr="00667113"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Ron was in charge of collecting the data. We started with data from Kaggle, but it was found out that it is not working. While Our project requires to use label data, the hashtags in Kaggle’s data were not significantly towards positive or negative directions. That is why we decided to scrape our one data, with hashtags we can label. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> In order to choose the right hashtags, w</w:t></w:r><w:r><w:t xml:space="preserve">e took the existing data set and viewed hashtags from the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>kaggle</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> set</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>. Then, we p</w:t></w:r><w:r><w:t>ick</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>ed</w:t></w:r><w:r><w:t xml:space="preserve"> ones that were being used in relation to the vaccine that were clearly pro or against</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">We used </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00667113"><w:rPr><w:b/><w:bCs/><w:lang w:val="en-US"/></w:rPr><w:t>tweepy</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00667113"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> for the scraping, together with Twitter developer account credentials. We faced scraping limitations made by Twitter, which made the process harder. Another obstacle we faced was that we can scrape data from tweepy only from the previous week.</w:t></w:r></w:p>
'@

$rng.InsertXML($newParagraphXml)
